$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in the title cell
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 21 de Marzo de 2020 a las 03:46"

# Row 8: Iran -> Estados Unidos
$ws.Cells.Item(8, 1).Value = "Estados Unidos"
$ws.Cells.Item(8, 2).Value = 19650
$ws.Cells.Item(8, 3).Value = 5861
$ws.Cells.Item(8, 4).Value = 147
$ws.Cells.Item(8, 5).Value = 19240
$ws.Cells.Item(8, 6).Value = 64
$ws.Cells.Item(8, 7).Value = 56
$ws.Cells.Item(8, 8).Value = 263

# Row 9: Estados Unidos -> Iran
$ws.Cells.Item(9, 1).Value = "Iran"
$ws.Cells.Item(9, 2).Value = 19644
$ws.Cells.Item(9, 3).Value = 1237
$ws.Cells.Item(9, 4).Value = 6745
$ws.Cells.Item(9, 5).Value = 11466
$ws.Cells.Item(9, 6).Value = 0
$ws.Cells.Item(9, 7).Value = 149
$ws.Cells.Item(9, 8).Value = 1433

# Row 12: Suiza -> Suiza
$ws.Cells.Item(12, 6).Value = 141

# Row 23: Brasil -> Japon
$ws.Cells.Item(23, 1).Value = "Japon"
$ws.Cells.Item(23, 2).Value = 1007
$ws.Cells.Item(23, 3).Value = 64
$ws.Cells.Item(23, 4).Value = 215
$ws.Cells.Item(23, 5).Value = 757
$ws.Cells.Item(23, 6).Value = 50
$ws.Cells.Item(23, 7).Value = 2
$ws.Cells.Item(23, 8).Value = 35

# Row 24: Japon -> Brasil
$ws.Cells.Item(24, 1).Value = "Brasil"
$ws.Cells.Item(24, 2).Value = 970
$ws.Cells.Item(24, 3).Value = 330
$ws.Cells.Item(24, 4).Value = 2
$ws.Cells.Item(24, 5).Value = 957
$ws.Cells.Item(24, 6).Value = 18
$ws.Cells.Item(24, 7).Value = 4
$ws.Cells.Item(24, 8).Value = 11

# Row 27: Crucero -> Crucero
$ws.Cells.Item(27, 4).Value = 567
$ws.Cells.Item(27, 5).Value = 137
$ws.Cells.Item(27, 6).Value = 15

# Row 96: Oman -> Guadalupe
$ws.Cells.Item(96, 1).Value = "Guadalupe"
$ws.Cells.Item(96, 2).Value = 51
$ws.Cells.Item(96, 3).Value = 18
$ws.Cells.Item(96, 4).Value = 0
$ws.Cells.Item(96, 5).Value = 50
$ws.Cells.Item(96, 6).Value = 4
$ws.Cells.Item(96, 7).Value = 1
$ws.Cells.Item(96, 8).Value = 1

# Row 97: Estado de Palestina -> Oman
$ws.Cells.Item(97, 1).Value = "Oman"
$ws.Cells.Item(97, 3).Value = 0
$ws.Cells.Item(97, 4).Value = 13
$ws.Cells.Item(97, 5).Value = 35

# Row 98: Senegal -> Estado de Palestina
$ws.Cells.Item(98, 1).Value = "Estado de Palestina"
$ws.Cells.Item(98, 2).Value = 48
$ws.Cells.Item(98, 3).Value = 1
$ws.Cells.Item(98, 4).Value = 17
$ws.Cells.Item(98, 5).Value = 31

# Row 99: Guadalupe -> Senegal
$ws.Cells.Item(99, 1).Value = "Senegal"
$ws.Cells.Item(99, 2).Value = 47
$ws.Cells.Item(99, 3).Value = 11
$ws.Cells.Item(99, 4).Value = 5
$ws.Cells.Item(99, 5).Value = 42

# Row 104: Uzbekistan -> Reunion
$ws.Cells.Item(104, 1).Value = "Reunion"
$ws.Cells.Item(104, 2).Value = 38
$ws.Cells.Item(104, 5).Value = 38

# Row 105: Martinica -> Uzbekistan
$ws.Cells.Item(105, 1).Value = "Uzbekistan"
$ws.Cells.Item(105, 2).Value = 33
$ws.Cells.Item(105, 3).Value = 10
$ws.Cells.Item(105, 5).Value = 33
$ws.Cells.Item(105, 6).Value = 0
$ws.Cells.Item(105, 8).Value = 0

# Row 106: Liechtenstein -> Martinica
$ws.Cells.Item(106, 1).Value = "Martinica"
$ws.Cells.Item(106, 2).Value = 32
$ws.Cells.Item(106, 3).Value = 9
$ws.Cells.Item(106, 5).Value = 31
$ws.Cells.Item(106, 6).Value = 7
$ws.Cells.Item(106, 8).Value = 1

# Row 107: Reunion -> Liechtenstein
$ws.Cells.Item(107, 1).Value = "Liechtenstein"

# Row 108: Honduras -> Camerun
$ws.Cells.Item(108, 1).Value = "Camerun"
$ws.Cells.Item(108, 2).Value = 27
$ws.Cells.Item(108, 3).Value = 14
$ws.Cells.Item(108, 4).Value = 2
$ws.Cells.Item(108, 5).Value = 25

# Row 109: Afganistan -> Honduras
$ws.Cells.Item(109, 1).Value = "Honduras"
$ws.Cells.Item(109, 3).Value = 12
$ws.Cells.Item(109, 4).Value = 0
$ws.Cells.Item(109, 5).Value = 24

# Row 110: Cuba -> Afganistan
$ws.Cells.Item(110, 1).Value = "Afganistan"
$ws.Cells.Item(110, 2).Value = 24
$ws.Cells.Item(110, 3).Value = 2
$ws.Cells.Item(110, 4).Value = 1
$ws.Cells.Item(110, 5).Value = 23
$ws.Cells.Item(110, 8).Value = 0

# Row 111: Camerun -> Cuba
$ws.Cells.Item(111, 1).Value = "Cuba"
$ws.Cells.Item(111, 2).Value = 21
$ws.Cells.Item(111, 3).Value = 10
$ws.Cells.Item(111, 4).Value = 0
$ws.Cells.Item(111, 5).Value = 20
$ws.Cells.Item(111, 8).Value = 1

# Row 114: Paraguay -> Consejo Danes para los Refugiados
$ws.Cells.Item(114, 1).Value = "Consejo Danes para los Refugiados"
$ws.Cells.Item(114, 3).Value = 4
$ws.Cells.Item(114, 6).Value = 0

# Row 115: Consejo Danes para los Refugiados -> Paraguay
$ws.Cells.Item(115, 1).Value = "Paraguay"
$ws.Cells.Item(115, 3).Value = 5
$ws.Cells.Item(115, 6).Value = 1

# Row 118: Ghana -> Bolivia
$ws.Cells.Item(118, 1).Value = "Bolivia"
$ws.Cells.Item(118, 3).Value = 1

# Row 119: Bolivia -> Ghana
$ws.Cells.Item(119, 1).Value = "Ghana"
$ws.Cells.Item(119, 3).Value = 5

# Row 122: Montenegro -> Guam
$ws.Cells.Item(122, 1).Value = "Guam"
$ws.Cells.Item(122, 3).Value = 2

# Row 123: Guam -> Montenegro
$ws.Cells.Item(123, 1).Value = "Montenegro"
$ws.Cells.Item(123, 3).Value = 1

# Row 132: Trinidad yTobago -> Togo
$ws.Cells.Item(132, 1).Value = "Togo"
$ws.Cells.Item(132, 3).Value = 8

# Row 133: Etiopia -> Trinidad yTobago
$ws.Cells.Item(133, 1).Value = "Trinidad yTobago"
$ws.Cells.Item(133, 3).Value = 0

# Row 134: Togo -> Etiopia
$ws.Cells.Item(134, 1).Value = "Etiopia"
$ws.Cells.Item(134, 3).Value = 2

# Row 136: Kenia -> Seychelles
$ws.Cells.Item(136, 1).Value = "Seychelles"
$ws.Cells.Item(136, 3).Value = 1

# Row 137: Seychelles -> Mayotte
$ws.Cells.Item(137, 1).Value = "Mayotte"
$ws.Cells.Item(137, 3).Value = 3

# Row 138: Kirguistan -> Kenia
$ws.Cells.Item(138, 1).Value = "Kenia"
$ws.Cells.Item(138, 2).Value = 7
$ws.Cells.Item(138, 3).Value = 0
$ws.Cells.Item(138, 5).Value = 7

# Row 139: Mayotte -> Kirguistan
$ws.Cells.Item(139, 1).Value = "Kirguistan"
$ws.Cells.Item(139, 3).Value = 3

# Row 163: Mauritania -> Haiti
$ws.Cells.Item(163, 1).Value = "Haiti"
$ws.Cells.Item(163, 3).Value = 2

# Row 166: Haiti -> Mauritania
$ws.Cells.Item(166, 1).Value = "Mauritania"
$ws.Cells.Item(166, 3).Value = 0
